# Rows 12-14 get their species-related data cyclically rotated:
#   new row 12 = old row 13 data
#   new row 13 = old row 14 data
#   new row 14 = old row 12 data
# Only columns A, B, E, F, G, H, P, Q, R, S are affected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "P", "Q", "R", "S")

# Capture the original values for rows 12, 13, 14 before overwriting anything.
$orig = @{}
foreach ($r in 12, 13, 14) {
    $orig[$r] = @{}
    foreach ($c in $cols) {
        $orig[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# New row 12 <- old row 13
foreach ($c in $cols) {
    $ws.Range("$c" + "12").Value = $orig[13][$c]
}

# New row 13 <- old row 14
foreach ($c in $cols) {
    $ws.Range("$c" + "13").Value = $orig[14][$c]
}

# New row 14 <- old row 12
foreach ($c in $cols) {
    $ws.Range("$c" + "14").Value = $orig[12][$c]
}
